$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to the "展览" and
# "全部类型" sheets (rows 2-18).
$updates = @{
    2  = 247
    3  = 274
    4  = 284
    5  = 827
    6  = 283
    7  = 6762
    8  = 57
    11 = 85
    12 = 39
    16 = 230
    17 = 581
    18 = 65
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
